$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 92513.55
$ws.Range("I6").Value = 92513.55
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 277540.65
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -277428.65
$ws.Range("H17").Value = 244223.36
$ws.Range("J17").Value = 250266.45
$ws.Range("L17").Value = 750799.3500000001
$ws.Range("N17").Value = -751135.3500000001
$ws.Range("H19").Value = 4198.125
$ws.Range("I19").Value = 2117.2
$ws.Range("K19").Value = 2117.2
$ws.Range("M19").Value = -1942.2
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("H28").Value = 4265.2144
$ws.Range("I28").Value = 4226.25
$ws.Range("K28").Value = 4226.25
$ws.Range("M28").Value = -3741.25
$ws.Range("H33").Value = 2475.5
$ws.Range("I33").Value = 649.1
$ws.Range("K33").Value = 649.1
$ws.Range("M33").Value = -420.1
$ws.Range("H41").Value = 717.2727
$ws.Range("J41").Value = 982.3333
$ws.Range("L41").Value = 982.3333
$ws.Range("N41").Value = -1862.3333
$ws.Range("H48").Value = 1788.2
$ws.Range("J48").Value = 1788.2
$ws.Range("L48").Value = 5364.6
$ws.Range("N48").Value = -5948.6
$ws.Range("H55").Value = 507
$ws.Range("I55").Value = 175
$ws.Range("J55").Value = 1088
$ws.Range("K55").Value = 175
$ws.Range("L55").Value = 1088
$ws.Range("M55").Value = 39
$ws.Range("N55").Value = -1516
$ws.Range("H56").Value = 1788.2
$ws.Range("J56").Value = 1788.2
$ws.Range("L56").Value = 5364.6
$ws.Range("N56").Value = -6432.6
$ws.Range("H69").Value = 9747.299999999999
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H70").Value = 5200.7046
$ws.Range("I70").Value = 1069.8235
$ws.Range("J70").Value = 7801.6294
$ws.Range("K70").Value = 3209.4705
$ws.Range("L70").Value = 23404.8882
$ws.Range("M70").Value = -2939.4705
$ws.Range("N70").Value = -23944.8882
$ws.Range("H72").Value = 9747.299999999999
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H73").Value = 5200.7046
$ws.Range("I73").Value = 1069.8235
$ws.Range("J73").Value = 7801.6294
$ws.Range("K73").Value = 3209.4705
$ws.Range("L73").Value = 23404.8882
$ws.Range("M73").Value = -2273.4705
$ws.Range("N73").Value = -25276.8882
$ws.Range("H98").Value = 1979983.2
$ws.Range("I98").Value = 2067469.6
$ws.Range("K98").Value = 2067469.6
$ws.Range("M98").Value = -2065971.6
$ws.Range("H101").Value = 3496825.2
$ws.Range("I101").Value = 9091332
$ws.Range("J101").Value = 258.25
$ws.Range("K101").Value = 27273996
$ws.Range("L101").Value = 774.75
$ws.Range("M101").Value = -27272374
$ws.Range("N101").Value = -4018.75
$ws.Range("H106").Value = 2182.9443
$ws.Range("I106").Value = 1819.7333
$ws.Range("K106").Value = 1819.7333
$ws.Range("M106").Value = -1188.7333
$ws.Range("H111").Value = 849.5
$ws.Range("I111").Value = 799.5454999999999
$ws.Range("K111").Value = 2398.6365
$ws.Range("M111").Value = 668.3635000000004
$ws.Range("H116").Value = 5119.5454
$ws.Range("I116").Value = 4404
$ws.Range("K116").Value = 4404
$ws.Range("M116").Value = -962
$ws.Range("H122").Value = 1979983.2
$ws.Range("I122").Value = 2067469.6
$ws.Range("K122").Value = 6202408.800000001
$ws.Range("M122").Value = -6199958.800000001
$ws.Range("H135").Value = 949.45
$ws.Range("I135").Value = 949.45
$ws.Range("K135").Value = 8545.050000000001
$ws.Range("M135").Value = -6010.050000000001
$ws.Range("H138").Value = 2704.528
$ws.Range("I138").Value = 1354.625
$ws.Range("K138").Value = 4063.875
$ws.Range("M138").Value = 1076.125
$ws.Range("N6").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("M23").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2998
$ws.Range("I45").Value = 2998
$ws.Range("K45").Value = 2998
$ws.Range("M45").Value = -2621
$ws.Range("H61").Value = 5172.2856
$ws.Range("I61").Value = 3182.9092
$ws.Range("K61").Value = 3182.9092
$ws.Range("M61").Value = -2970.9092
$ws.Range("H74").Value = 111530.89
$ws.Range("I74").Value = 132769.14
$ws.Range("J74").Value = 5339.6665
$ws.Range("K74").Value = 132769.14
$ws.Range("L74").Value = 5339.6665
$ws.Range("M74").Value = -131895.14
$ws.Range("N74").Value = -7087.6665
$ws.Range("H77").Value = 111530.89
$ws.Range("I77").Value = 132769.14
$ws.Range("J77").Value = 5339.6665
$ws.Range("K77").Value = 663845.7000000001
$ws.Range("L77").Value = 26698.3325
$ws.Range("M77").Value = -659477.7000000001
$ws.Range("N77").Value = -35434.3325
$ws.Range("H110").Value = 2693.44
$ws.Range("I110").Value = 1428.9546
$ws.Range("K110").Value = 1428.9546
$ws.Range("M110").Value = 616.0454
$ws.Range("H122").Value = 5949.5
$ws.Range("I122").Value = 2932.6667
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 8798.000100000001
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -6348.000100000001
$ws.Range("N122").Value = -49900
$ws.Range("H132").Value = 3490
$ws.Range("I132").Value = 3301.375
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9904.125
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -7374.125
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 5172.2856
$ws.Range("I136").Value = 3182.9092
$ws.Range("K136").Value = 9548.7276
$ws.Range("M136").Value = -6998.7276

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2186.5
$ws.Range("J20").Value = 2680
$ws.Range("L20").Value = 2680
$ws.Range("N20").Value = -3174
$ws.Range("H86").Value = 5547.1
$ws.Range("I86").Value = 5746.4287
$ws.Range("J86").Value = 5082
$ws.Range("K86").Value = 5746.4287
$ws.Range("L86").Value = 5082
$ws.Range("M86").Value = -4623.4287
$ws.Range("N86").Value = -7328
$ws.Range("H89").Value = 5547.1
$ws.Range("I89").Value = 5746.4287
$ws.Range("J89").Value = 5082
$ws.Range("K89").Value = 28732.1435
$ws.Range("L89").Value = 25410
$ws.Range("M89").Value = -23116.1435
$ws.Range("N89").Value = -36642
$ws.Range("H105").Value = 2071.9678
$ws.Range("I105").Value = 2053.5186
$ws.Range("J105").Value = 2196.5
$ws.Range("K105").Value = 2053.5186
$ws.Range("L105").Value = 2196.5
$ws.Range("M105").Value = -306.5185999999999
$ws.Range("N105").Value = -5690.5
$ws.Range("H134").Value = 3931.611
$ws.Range("I134").Value = 3105.6924
$ws.Range("J134").Value = 6079
$ws.Range("K134").Value = 9317.0772
$ws.Range("L134").Value = 18237
$ws.Range("M134").Value = -6782.0772
$ws.Range("N134").Value = -23307

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2918
$ws.Range("I16").Value = 2987.2
$ws.Range("J16").Value = 2745
$ws.Range("K16").Value = 2987.2
$ws.Range("L16").Value = 2745
$ws.Range("M16").Value = -2700.2
$ws.Range("N16").Value = -3319
$ws.Range("H58").Value = 6251.6665
$ws.Range("I58").Value = 3177.889
$ws.Range("K58").Value = 3177.889
$ws.Range("M58").Value = -2974.889
$ws.Range("H64").Value = 41421.082
$ws.Range("J64").Value = 41421.082
$ws.Range("L64").Value = 41421.082
$ws.Range("N64").Value = -41917.082
$ws.Range("H67").Value = 41421.082
$ws.Range("J67").Value = 41421.082
$ws.Range("L67").Value = 41421.082
$ws.Range("N67").Value = -43137.082
$ws.Range("H97").Value = 50196.25
$ws.Range("J97").Value = 50193.5
$ws.Range("L97").Value = 50193.5
$ws.Range("N97").Value = -52175.5
$ws.Range("H99").Value = 335535.9
$ws.Range("I99").Value = 503802.66
$ws.Range("J99").Value = 29596.363
$ws.Range("K99").Value = 503802.66
$ws.Range("L99").Value = 29596.363
$ws.Range("M99").Value = -502304.66
$ws.Range("N99").Value = -32592.363
$ws.Range("H109").Value = 35284.5
$ws.Range("J109").Value = 35284.5
$ws.Range("L109").Value = 35284.5
$ws.Range("N109").Value = -37364.5
$ws.Range("H113").Value = 2918
$ws.Range("I113").Value = 2987.2
$ws.Range("J113").Value = 2745
$ws.Range("K113").Value = 2987.2
$ws.Range("L113").Value = 2745
$ws.Range("M113").Value = -817.1999999999998
$ws.Range("N113").Value = -7085
$ws.Range("H122").Value = 2128.375
$ws.Range("I122").Value = 1811.5
$ws.Range("J122").Value = 2234
$ws.Range("K122").Value = 5434.5
$ws.Range("L122").Value = 6702
$ws.Range("M122").Value = -2984.5
$ws.Range("N122").Value = -11602
$ws.Range("H126").Value = 335535.9
$ws.Range("I126").Value = 503802.66
$ws.Range("J126").Value = 29596.363
$ws.Range("K126").Value = 1511407.98
$ws.Range("L126").Value = 88789.08900000001
$ws.Range("M126").Value = -1508937.98
$ws.Range("N126").Value = -93729.08900000001
$ws.Range("H132").Value = 2831.5806
$ws.Range("I132").Value = 2129.5217
$ws.Range("K132").Value = 6388.5651
$ws.Range("M132").Value = -3858.5651
$ws.Range("H134").Value = 8088.4287
$ws.Range("I134").Value = 9432.25
$ws.Range("K134").Value = 28296.75
$ws.Range("M134").Value = -25761.75
$ws.Range("H136").Value = 6251.6665
$ws.Range("I136").Value = 3177.889
$ws.Range("K136").Value = 9533.667000000001
$ws.Range("M136").Value = -6983.667000000001
$ws.Range("H140").Value = 69166.664
$ws.Range("J140").Value = 69166.664
$ws.Range("L140").Value = 69166.664
$ws.Range("N140").Value = -79526.664

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 133727.53
$ws.Range("I11").Value = 167151.5
$ws.Range("J11").Value = 31.666666
$ws.Range("K11").Value = 501454.5
$ws.Range("L11").Value = 94.99999800000001
$ws.Range("M11").Value = -501314.5
$ws.Range("N11").Value = -374.999998
$ws.Range("H34").Value = 838.5172
$ws.Range("I34").Value = 787.5
$ws.Range("J34").Value = 1280.6666
$ws.Range("K34").Value = 2362.5
$ws.Range("L34").Value = 3841.9998
$ws.Range("M34").Value = -2278.5
$ws.Range("N34").Value = -4009.9998
$ws.Range("H39").Value = 949.0454999999999
$ws.Range("J39").Value = 2707.1667
$ws.Range("L39").Value = 8121.500100000001
$ws.Range("N39").Value = -8709.500100000001
$ws.Range("H55").Value = 5841.4375
$ws.Range("J55").Value = 9042.857
$ws.Range("L55").Value = 27128.571
$ws.Range("N55").Value = -27482.571
$ws.Range("H68").Value = 16669147
$ws.Range("J68").Value = 9999
$ws.Range("L68").Value = 29997
$ws.Range("N68").Value = -31619
$ws.Range("H71").Value = 16669147
$ws.Range("J71").Value = 9999
$ws.Range("L71").Value = 89991
$ws.Range("N71").Value = -98103
$ws.Range("H86").Value = 865.6667
$ws.Range("I86").Value = 878.5714
$ws.Range("J86").Value = 854.375
$ws.Range("K86").Value = 2635.7142
$ws.Range("L86").Value = 2563.125
$ws.Range("M86").Value = -1449.7142
$ws.Range("N86").Value = -4935.125
$ws.Range("H89").Value = 865.6667
$ws.Range("I89").Value = 878.5714
$ws.Range("J89").Value = 854.375
$ws.Range("K89").Value = 7907.1426
$ws.Range("L89").Value = 7689.375
$ws.Range("M89").Value = -1979.1426
$ws.Range("N89").Value = -19545.375
$ws.Range("H103").Value = 1033
$ws.Range("J103").Value = 1106
$ws.Range("L103").Value = 3318
$ws.Range("N103").Value = -5076
$ws.Range("H117").Value = 392.6
$ws.Range("I117").Value = 240.75
$ws.Range("K117").Value = 722.25
$ws.Range("M117").Value = 2719.75
$ws.Range("H121").Value = 96325.27
$ws.Range("J121").Value = 116398
$ws.Range("L121").Value = 349194
$ws.Range("N121").Value = -351814
$ws.Range("H129").Value = 20834978
$ws.Range("I129").Value = 41668100
$ws.Range("J129").Value = 1856.625
$ws.Range("K129").Value = 125004300
$ws.Range("L129").Value = 5569.875
$ws.Range("M129").Value = -124999300
$ws.Range("N129").Value = -15569.875
$ws.Range("H131").Value = 7144224.5
$ws.Range("J131").Value = 1412.3077
$ws.Range("L131").Value = 4236.9231
$ws.Range("N131").Value = -14316.9231
$ws.Range("H132").Value = 4606.575
$ws.Range("J132").Value = 2498.375
$ws.Range("L132").Value = 22485.375
$ws.Range("N132").Value = -27545.375
$ws.Range("H138").Value = 2344.2856
$ws.Range("I138").Value = 2344.2856
$ws.Range("K138").Value = 7032.8568
$ws.Range("M138").Value = -1892.8568
$ws.Range("H139").Value = 2541.5908
$ws.Range("I139").Value = 1767.6875
$ws.Range("K139").Value = 5303.0625
$ws.Range("M139").Value = -163.0625
$ws.Range("H140").Value = 53160
$ws.Range("I140").Value = 100765
$ws.Range("K140").Value = 302295
$ws.Range("M140").Value = -297115

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 17998.2
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 17998.2
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 17998.2
$ws.Range("N20").Value = -18488.2
$ws.Range("H24").Value = 13992
$ws.Range("J24").Value = 13992
$ws.Range("L24").Value = 13992
$ws.Range("N24").Value = -14338
$ws.Range("H80").Value = 4165.4443
$ws.Range("J80").Value = 4315.931
$ws.Range("L80").Value = 4315.931
$ws.Range("N80").Value = -6311.931
$ws.Range("H83").Value = 4165.4443
$ws.Range("J83").Value = 4315.931
$ws.Range("L83").Value = 21579.655
$ws.Range("N83").Value = -31563.655
$ws.Range("H102").Value = 8406.6
$ws.Range("I102").Value = 8406.6
$ws.Range("K102").Value = 8406.6
$ws.Range("M102").Value = -6784.6
$ws.Range("H113").Value = 5634.1333
$ws.Range("I113").Value = 6658.1875
$ws.Range("K113").Value = 6658.1875
$ws.Range("M113").Value = -4488.1875
$ws.Range("H122").Value = 4508.4614
$ws.Range("I122").Value = 8196.666999999999
$ws.Range("J122").Value = 1347.1428
$ws.Range("K122").Value = 24590.001
$ws.Range("L122").Value = 4041.4284
$ws.Range("M122").Value = -22140.001
$ws.Range("N122").Value = -8941.428400000001
$ws.Range("H126").Value = 6220.1333
$ws.Range("I126").Value = 7725.5557
$ws.Range("J126").Value = 3962
$ws.Range("K126").Value = 23176.6671
$ws.Range("L126").Value = 11886
$ws.Range("M126").Value = -20706.6671
$ws.Range("N126").Value = -16826
$ws.Range("H127").Value = 46000
$ws.Range("J127").Value = 46000
$ws.Range("L127").Value = 46000
$ws.Range("N127").Value = -55920
$ws.Range("H132").Value = 1000012
$ws.Range("I132").Value = 1000012
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000036
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2997506
$ws.Range("M20").ClearContents()
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4603.815
$ws.Range("I7").Value = 5164.9565
$ws.Range("J7").Value = 1377.25
$ws.Range("K7").Value = 5164.9565
$ws.Range("L7").Value = 1377.25
$ws.Range("M7").Value = -5052.9565
$ws.Range("N7").Value = -1601.25
$ws.Range("H16").Value = 843.26666
$ws.Range("I16").Value = 948.53845
$ws.Range("K16").Value = 948.53845
$ws.Range("M16").Value = -778.53845
$ws.Range("H20").Value = 7800.6
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 12962.333
$ws.Range("K20").Value = 58
$ws.Range("L20").Value = 12962.333
$ws.Range("M20").Value = 168
$ws.Range("N20").Value = -13414.333
$ws.Range("H22").Value = 1170.8064
$ws.Range("J22").Value = 1332.1
$ws.Range("L22").Value = 1332.1
$ws.Range("N22").Value = -1922.1
$ws.Range("H27").Value = 1170.8064
$ws.Range("J27").Value = 1332.1
$ws.Range("L27").Value = 1332.1
$ws.Range("N27").Value = -1546.1
$ws.Range("H29").Value = 29999
$ws.Range("I29").Value = 29999
$ws.Range("K29").Value = 29999
$ws.Range("M29").Value = -29704
$ws.Range("H40").Value = 14911.333
$ws.Range("I40").Value = 25249
$ws.Range("J40").Value = 4573.6665
$ws.Range("K40").Value = 25249
$ws.Range("L40").Value = 4573.6665
$ws.Range("M40").Value = -25113
$ws.Range("N40").Value = -4845.6665
$ws.Range("H46").Value = 3800.0908
$ws.Range("J46").Value = 4667.6665
$ws.Range("L46").Value = 4667.6665
$ws.Range("N46").Value = -5043.6665
$ws.Range("H82").Value = 1544.5
$ws.Range("I82").Value = 1216.9
$ws.Range("J82").Value = 1954
$ws.Range("K82").Value = 1216.9
$ws.Range("L82").Value = 1954
$ws.Range("M82").Value = -855.9000000000001
$ws.Range("N82").Value = -2676
$ws.Range("H85").Value = 1544.5
$ws.Range("I85").Value = 1216.9
$ws.Range("J85").Value = 1954
$ws.Range("K85").Value = 1216.9
$ws.Range("L85").Value = 1954
$ws.Range("M85").Value = 31.09999999999991
$ws.Range("N85").Value = -4450
$ws.Range("H126").Value = 4603.815
$ws.Range("I126").Value = 5164.9565
$ws.Range("J126").Value = 1377.25
$ws.Range("K126").Value = 15494.8695
$ws.Range("L126").Value = 4131.75
$ws.Range("M126").Value = -13024.8695
$ws.Range("N126").Value = -9071.75
$ws.Range("H132").Value = 12434.454
$ws.Range("I132").Value = 13467.9
$ws.Range("K132").Value = 40403.7
$ws.Range("M132").Value = -37873.7
$ws.Range("H136").Value = 1339.8667
$ws.Range("I136").Value = 1249.8572
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 3749.5716
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -1199.5716
$ws.Range("N136").Value = -12900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 9332.75
$ws.Range("J31").Value = 10777
$ws.Range("L31").Value = 10777
$ws.Range("N31").Value = -11473
$ws.Range("H34").Value = 29757.25
$ws.Range("I34").Value = 29757.25
$ws.Range("K34").Value = 29757.25
$ws.Range("M34").Value = -29554.25
$ws.Range("H51").Value = 15009.533
$ws.Range("I51").Value = 18295.572
$ws.Range("J51").Value = 12134.25
$ws.Range("K51").Value = 18295.572
$ws.Range("L51").Value = 12134.25
$ws.Range("M51").Value = -17785.572
$ws.Range("N51").Value = -13154.25
$ws.Range("H52").Value = 13181.333
$ws.Range("J52").Value = 22044.5
$ws.Range("L52").Value = 22044.5
$ws.Range("N52").Value = -22496.5
$ws.Range("H62").Value = 5917.467
$ws.Range("I62").Value = 2522.75
$ws.Range("J62").Value = 9797.143
$ws.Range("K62").Value = 2522.75
$ws.Range("L62").Value = 9797.143
$ws.Range("M62").Value = -1898.75
$ws.Range("N62").Value = -11045.143
$ws.Range("H65").Value = 5917.467
$ws.Range("I65").Value = 2522.75
$ws.Range("J65").Value = 9797.143
$ws.Range("K65").Value = 12613.75
$ws.Range("L65").Value = 48985.715
$ws.Range("M65").Value = -9493.75
$ws.Range("N65").Value = -55225.715
$ws.Range("H113").Value = 631.1
$ws.Range("I113").Value = 567.8889
$ws.Range("K113").Value = 1703.6667
$ws.Range("M113").Value = 466.3332999999998
$ws.Range("H122").Value = 10151.243
$ws.Range("I122").Value = 2700.3438
$ws.Range("J122").Value = 57837
$ws.Range("K122").Value = 8101.0314
$ws.Range("L122").Value = 173511
$ws.Range("M122").Value = -5651.0314
$ws.Range("N122").Value = -178411
$ws.Range("H126").Value = 2031.7368
$ws.Range("I126").Value = 2020.6129
$ws.Range("K126").Value = 6061.8387
$ws.Range("M126").Value = -3591.8387
$ws.Range("H132").Value = 1704.175
$ws.Range("I132").Value = 1687.8918
$ws.Range("K132").Value = 5063.6754
$ws.Range("M132").Value = -2533.6754
$ws.Range("H136").Value = 810280.7
$ws.Range("I136").Value = 836137.4399999999
$ws.Range("K136").Value = 2508412.32
$ws.Range("M136").Value = -2505862.32
